$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2848.5
$ws.Range("I70").Value = 1533.1111
$ws.Range("K70").Value = 4599.3333
$ws.Range("M70").Value = -4329.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 2848.5
$ws.Range("I73").Value = 1533.1111
$ws.Range("K73").Value = 4599.3333
$ws.Range("M73").Value = -3663.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 5278.517
$ws.Range("I86").Value = 10300.909
$ws.Range("K86").Value = 10300.909
$ws.Range("M86").Value = -9177.909

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 5278.517
$ws.Range("I89").Value = 10300.909
$ws.Range("K89").Value = 51504.545
$ws.Range("M89").Value = -45888.545

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2760.5386
$ws.Range("I98").Value = 1860.4286
$ws.Range("J98").Value = 6541
$ws.Range("K98").Value = 1860.4286
$ws.Range("L98").Value = 6541
$ws.Range("M98").Value = -362.4286
$ws.Range("N98").Value = -9537

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 1285.4615
$ws.Range("I106").Value = 1211.6666
$ws.Range("J106").Value = 1451.5
$ws.Range("K106").Value = 1211.6666
$ws.Range("L106").Value = 1451.5
$ws.Range("M106").Value = -580.6666
$ws.Range("N106").Value = -2713.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 3452.6667
$ws.Range("I118").Value = 566.25
$ws.Range("J118").Value = 5761.8
$ws.Range("K118").Value = 1698.75
$ws.Range("L118").Value = 17285.4
$ws.Range("M118").Value = -41.75
$ws.Range("N118").Value = -20599.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 2760.5386
$ws.Range("I122").Value = 1860.4286
$ws.Range("J122").Value = 6541
$ws.Range("K122").Value = 5581.2858
$ws.Range("L122").Value = 19623
$ws.Range("M122").Value = -3131.2858
$ws.Range("N122").Value = -24523

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 4099915.2
$ws.Range("J129").Value = 1586.2181
$ws.Range("L129").Value = 4758.6543
$ws.Range("N129").Value = -14758.6543

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 6255637.5
$ws.Range("I137").Value = 25015500
$ws.Range("J137").Value = 2350.25
$ws.Range("K137").Value = 75046500
$ws.Range("L137").Value = 7050.75
$ws.Range("M137").Value = -75043950
$ws.Range("N137").Value = -12150.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 12502629
$ws.Range("I2").Value = 19232314
$ws.Range("J2").Value = 4642.857
$ws.Range("K2").Value = 19232314
$ws.Range("L2").Value = 4642.857
$ws.Range("M2").Value = -19232201
$ws.Range("N2").Value = -4868.857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3872.2
$ws.Range("I61").Value = 1871.8572
$ws.Range("J61").Value = 5622.5
$ws.Range("K61").Value = 1871.8572
$ws.Range("L61").Value = 5622.5
$ws.Range("M61").Value = -1659.8572
$ws.Range("N61").Value = -6046.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 12502629
$ws.Range("I116").Value = 19232314
$ws.Range("J116").Value = 4642.857
$ws.Range("K116").Value = 19232314
$ws.Range("L116").Value = 4642.857
$ws.Range("M116").Value = -19230020
$ws.Range("N116").Value = -9230.857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1059847.2
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 1456915
$ws.Range("K122").Value = 3000
$ws.Range("L122").Value = 4370745
$ws.Range("M122").Value = -550
$ws.Range("N122").Value = -4375645

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3872.2
$ws.Range("I136").Value = 1871.8572
$ws.Range("J136").Value = 5622.5
$ws.Range("K136").Value = 5615.571599999999
$ws.Range("L136").Value = 16867.5
$ws.Range("M136").Value = -3065.571599999999
$ws.Range("N136").Value = -21967.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 12502629
$ws.Range("I3").Value = 19232314
$ws.Range("J3").Value = 4642.857
$ws.Range("K3").Value = 19232314
$ws.Range("L3").Value = 4642.857
$ws.Range("M3").Value = -19232200
$ws.Range("N3").Value = -4870.857

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2313.158
$ws.Range("I20").Value = 1520.7273
$ws.Range("K20").Value = 1520.7273
$ws.Range("M20").Value = -1273.7273

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2922.2222
$ws.Range("I86").Value = 1716.6666
$ws.Range("J86").Value = 5333.3335
$ws.Range("K86").Value = 1716.6666
$ws.Range("L86").Value = 5333.3335
$ws.Range("M86").Value = -593.6666
$ws.Range("N86").Value = -7579.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2922.2222
$ws.Range("I89").Value = 1716.6666
$ws.Range("J89").Value = 5333.3335
$ws.Range("K89").Value = 8583.333000000001
$ws.Range("L89").Value = 26666.6675
$ws.Range("M89").Value = -2967.333000000001
$ws.Range("N89").Value = -37898.6675

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1792.3125
$ws.Range("I105").Value = 1191.2142
$ws.Range("J105").Value = 6000
$ws.Range("K105").Value = 1191.2142
$ws.Range("L105").Value = 6000
$ws.Range("M105").Value = 555.7858000000001
$ws.Range("N105").Value = -9494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2149.353
$ws.Range("I134").Value = 1721.1875
$ws.Range("K134").Value = 5163.5625
$ws.Range("M134").Value = -2628.5625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 26190.834
$ws.Range("I121").Value = 343.22223
$ws.Range("J121").Value = 103733.664
$ws.Range("K121").Value = 1029.66669
$ws.Range("L121").Value = 311200.992
$ws.Range("M121").Value = 280.33331
$ws.Range("N121").Value = -313820.992

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 70007
$ws.Range("I21").Value = 70007
$ws.Range("J21").Value = 70007
$ws.Range("K21").Value = 70007
$ws.Range("L21").Value = 70007
$ws.Range("M21").Value = -69834
$ws.Range("N21").Value = -70353

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H30").Value = 70007
$ws.Range("I30").Value = 70007
$ws.Range("J30").Value = 70007
$ws.Range("K30").Value = 70007
$ws.Range("L30").Value = 70007
$ws.Range("M30").Value = -69902
$ws.Range("N30").Value = -70217

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2622.111
$ws.Range("I7").Value = 1666.3334
$ws.Range("J7").Value = 3100
$ws.Range("K7").Value = 1666.3334
$ws.Range("L7").Value = 3100
$ws.Range("M7").Value = -1554.3334
$ws.Range("N7").Value = -3324

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2622.111
$ws.Range("I126").Value = 1666.3334
$ws.Range("J126").Value = 3100
$ws.Range("K126").Value = 4999.0002
$ws.Range("L126").Value = 9300
$ws.Range("M126").Value = -2529.0002
$ws.Range("N126").Value = -14240

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 70007
$ws.Range("J15").Value = 70007
$ws.Range("L15").Value = 70007
$ws.Range("N15").Value = -70583

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3036.5715
$ws.Range("I107").Value = 265
$ws.Range("J107").Value = 5115.25
$ws.Range("K107").Value = 795
$ws.Range("L107").Value = 15345.75
$ws.Range("M107").Value = 1125
$ws.Range("N107").Value = -19185.75
